$d = $word.ActiveDocument
$rng = $d.Content
$rng.Find.Execute("DD", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$rng.Collapse(0)  # wdCollapseEnd
$insStart = $rng.Start
$rng.InsertAfter("ddddd")
$insEnd = $insStart + 5
$newRng = $d.Range($insStart, $insEnd)
Write-Output ("newRng: [" + $newRng.Text + "] start=" + $newRng.Start + " end=" + $newRng.End)
$newRng.Bold = 1
$newRng.Bold = 0
Write-Output $d.Content.Text
